{"js": "const replacements = [\n  [\"85\u00d746=\", \"87\u00d726=\"],\n  [\"59\u00d798=\", \"22\u00d787=\"],\n  [\"42\u00d748=\", \"85\u00d755=\"],\n  [\"17\u00d757=\", \"99\u00d754=\"],\n  [\"11\u00d765=\", \"70\u00d765=\"],\n  [\"93\u00d715=\", \"95\u00d748=\"],\n  [\"49\u00d755=\", \"80\u00d780=\"],\n  [\"60\u00d760=\", \"81\u00d736=\"],\n  [\"80\u00d726=\", \"86\u00d773=\"],\n  [\"24\u00d774=\", \"46\u00d724=\"],\n  [\"99\u00d733=\", \"48\u00d749=\"],\n  [\"48\u00d776=\", \"50\u00d754=\"],\n  [\"41\u00d774=\", \"56\u00d732=\"],\n  [\"18\u00d747=\", \"51\u00d721=\"],\n  [\"34\u00d786=\", \"35\u00d785=\"],\n  [\"90\u00d791=\", \"41\u00d753=\"],\n  [\"53\u00d722=\", \"43\u00d750=\"],\n  [\"76\u00d781=\", \"43\u00d748=\"],\n  [\"48\u00d733=\", \"27\u00d770=\"],\n  [\"68\u00d785=\", \"59\u00d781=\"],\n  [\"99\u00d793=\", \"42\u00d768=\"],\n  [\"38\u00d795=\", \"24\u00d739=\"],\n  [\"58\u00d742=\", \"98\u00d723=\"],\n  [\"63\u00d767=\", \"78\u00d789=\"],\n  [\"66\u00d799=\", \"74\u00d755=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const searchResults = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  searchResults.load(\"text\");\n  await context.sync();\n\n  if (searchResults.items.length === 0) {\n    throw new Error(`Could not find text: ${oldText}`);\n  }\n\n  for (const range of searchResults.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old=\"85\u00d746=\"; new=\"87\u00d726=\"},\n    @{old=\"59\u00d798=\"; new=\"22\u00d787=\"},\n    @{old=\"42\u00d748=\"; new=\"85\u00d755=\"},\n    @{old=\"17\u00d757=\"; new=\"99\u00d754=\"},\n    @{old=\"11\u00d765=\"; new=\"70\u00d765=\"},\n    @{old=\"93\u00d715=\"; new=\"95\u00d748=\"},\n    @{old=\"49\u00d755=\"; new=\"80\u00d780=\"},\n    @{old=\"60\u00d760=\"; new=\"81\u00d736=\"},\n    @{old=\"80\u00d726=\"; new=\"86\u00d773=\"},\n    @{old=\"24\u00d774=\"; new=\"46\u00d724=\"},\n    @{old=\"99\u00d733=\"; new=\"48\u00d749=\"},\n    @{old=\"48\u00d776=\"; new=\"50\u00d754=\"},\n    @{old=\"41\u00d774=\"; new=\"56\u00d732=\"},\n    @{old=\"18\u00d747=\"; new=\"51\u00d721=\"},\n    @{old=\"34\u00d786=\"; new=\"35\u00d785=\"},\n    @{old=\"90\u00d791=\"; new=\"41\u00d753=\"},\n    @{old=\"53\u00d722=\"; new=\"43\u00d750=\"},\n    @{old=\"76\u00d781=\"; new=\"43\u00d748=\"},\n    @{old=\"48\u00d733=\"; new=\"27\u00d770=\"},\n    @{old=\"68\u00d785=\"; new=\"59\u00d781=\"},\n    @{old=\"99\u00d793=\"; new=\"42\u00d768=\"},\n    @{old=\"38\u00d795=\"; new=\"24\u00d739=\"},\n    @{old=\"58\u00d742=\"; new=\"98\u00d723=\"},\n    @{old=\"63\u00d767=\"; new=\"78\u00d789=\"},\n    @{old=\"66\u00d799=\"; new=\"74\u00d755=\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $r.old,   # FindText\n        $false,   # MatchCase\n        $false,   # MatchWholeWord\n        $false,   # MatchWildcards\n        $false,   # MatchSoundsLike\n        $false,   # MatchAllWordForms\n        $true,    # Forward\n        1,        # Wrap (wdFindContinue)\n        $false,   # Format\n        $r.new,   # ReplaceWith\n        2         # Replace (wdReplaceAll)\n    )\n}\n"}
